$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before DF (01-nov) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("DF1").EntireColumn.Insert()

$ws1.Range("DF1").Value = "01-nov"

for ($row = 2; $row -le 25; $row++) {
    $ws1.Cells.Item($row, 110).Value = "-"
}

# --- Sheet "Gaz": append new row for 2025-10-30 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A138").NumberFormat = "@"
$ws2.Range("A138").Value = "2025-10-30"
$ws2.Range("A138").Style = "Normal"
$ws2.Range("B138").Value = 29.8

# --- Sheet "CO2": append new row for 2025-10-30 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A138").NumberFormat = "@"
$ws3.Range("A138").Value = "2025-10-30"
$ws3.Range("A138").Style = "Normal"
$ws3.Range("B138").Value = 78.36
